$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values. Regenerate them per the new save_data values.
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 3
